$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("Q2").Value = 2.03
$ws.Range("R2").Value = 1.83

# Row 3
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 2
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 2.75
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("W3").Value = 9
$ws.Range("X3").Value = 19
$ws.Range("AB3").Value = 51
$ws.Range("AC3").Value = 7
$ws.Range("AE3").Value = 19
$ws.Range("AG3").Value = 501
$ws.Range("AH3").Value = 5.5
$ws.Range("AJ3").Value = 9.5
$ws.Range("AK3").Value = 17
$ws.Range("AL3").Value = 19
$ws.Range("AM3").Value = 41
$ws.Range("AO3").Value = 26
$ws.Range("AP3").Value = 41
$ws.Range("AQ3").Value = 101
$ws.Range("AR3").Value = 151
$ws.Range("AS3").Value = 401
$ws.Range("AT3").Value = 2.38
$ws.Range("AU3").Value = 9
$ws.Range("AX3").Value = 11
$ws.Range("AY3").Value = 26
$ws.Range("BA3").Value = 67
$ws.Range("BB3").Value = 251

# Row 4
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 3.4
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1.85

# Row 5
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 3.5
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.93

# Row 6
$ws.Range("G6").Value = 1.7
$ws.Range("H6").Value = 3.7

# Row 7
$ws.Range("J7").Value = 5
$ws.Range("U7").Value = 1.8
$ws.Range("V7").Value = 1.91
$ws.Range("AG7").Value = 251
$ws.Range("AI7").Value = 8
$ws.Range("AM7").Value = 26

# Row 8
$ws.Range("G8").Value = 2.15
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 3.25
$ws.Range("Q8").Value = 1.85
$ws.Range("R8").Value = 2
$ws.Range("Z8").Value = 21
$ws.Range("AB8").Value = 26
$ws.Range("AG8").Value = 201
$ws.Range("AI8").Value = 17
$ws.Range("AO8").Value = 12
$ws.Range("AX8").Value = 17
$ws.Range("AY8").Value = 23

# Row 9
$ws.Range("G9").Value = 2.38
$ws.Range("H9").Value = 3.5
$ws.Range("I9").Value = 2.7
$ws.Range("J9").Value = 3
$ws.Range("Q9").Value = 1.67
$ws.Range("R9").Value = 2.15
$ws.Range("S9").Value = 1.33
$ws.Range("T9").Value = 3.25
$ws.Range("U9").Value = 1.57
$ws.Range("V9").Value = 2.25
$ws.Range("AH9").Value = 11
$ws.Range("AI9").Value = 15
$ws.Range("AJ9").Value = 10
$ws.Range("AO9").Value = 13
$ws.Range("AT9").Value = 3.25

# Row 10
$ws.Range("G10").Value = 4.1
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 1.85
$ws.Range("J10").Value = 4.5
$ws.Range("L10").Value = 2.6
$ws.Range("Q10").Value = 2.08
$ws.Range("R10").Value = 1.73
$ws.Range("AC10").Value = 9
$ws.Range("AJ10").Value = 9
$ws.Range("AL10").Value = 17

# Row 11
$ws.Range("G11").Value = 1.7
$ws.Range("H11").Value = 3.9
$ws.Range("I11").Value = 4.33
$ws.Range("J11").Value = 2.3
$ws.Range("K11").Value = 2.3
$ws.Range("L11").Value = 4.75
$ws.Range("X11").Value = 9
$ws.Range("AF11").Value = 41
$ws.Range("AI11").Value = 23
$ws.Range("AL11").Value = 34
$ws.Range("AM11").Value = 34
$ws.Range("AS11").Value = 101
$ws.Range("AX11").Value = 23
$ws.Range("BA11").Value = 81
$ws.Range("BB11").Value = 151

# Row 12
$ws.Range("G12").Value = 2.8
$ws.Range("I12").Value = 2.5
$ws.Range("J12").Value = 3.6
$ws.Range("L12").Value = 3.25
$ws.Range("W12").Value = 8
$ws.Range("X12").Value = 13
$ws.Range("Y12").Value = 11
$ws.Range("Z12").Value = 29
$ws.Range("AA12").Value = 26
$ws.Range("AC12").Value = 8
$ws.Range("AH12").Value = 7.5
$ws.Range("AI12").Value = 12
$ws.Range("AN12").Value = 4.75
$ws.Range("AS12").Value = 201
$ws.Range("AV12").Value = 51
$ws.Range("AW12").Value = 4.5
$ws.Range("AX12").Value = 15
$ws.Range("BA12").Value = 81

# Row 13
$ws.Range("G13").Value = 1.38
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 7
$ws.Range("K13").Value = 2.75
$ws.Range("U13").Value = 1.67
$ws.Range("V13").Value = 2.1
$ws.Range("Z13").Value = 10
$ws.Range("AB13").Value = 21
$ws.Range("AJ13").Value = 21
$ws.Range("AK13").Value = 81
$ws.Range("AW13").Value = 8.5
$ws.Range("AX13").Value = 34

# Row 15
$ws.Range("G15").Value = 2.55
$ws.Range("I15").Value = 2.7
$ws.Range("J15").Value = 3.2
$ws.Range("L15").Value = 3.4
$ws.Range("W15").Value = 9
$ws.Range("X15").Value = 13
$ws.Range("Y15").Value = 10
$ws.Range("Z15").Value = 26
$ws.Range("AH15").Value = 9
$ws.Range("AI15").Value = 13
$ws.Range("AJ15").Value = 10
$ws.Range("AK15").Value = 29
$ws.Range("AL15").Value = 21
$ws.Range("AN15").Value = 4.5
$ws.Range("AQ15").Value = 41
$ws.Range("AW15").Value = 4.75
$ws.Range("AX15").Value = 15
$ws.Range("AZ15").Value = 51

# Row 16
$ws.Range("G16").Value = 2.3
$ws.Range("H16").Value = 3.25
$ws.Range("I16").Value = 3.1
$ws.Range("J16").Value = 3
$ws.Range("L16").Value = 3.75
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 10
$ws.Range("Q16").Value = 2.05
$ws.Range("R16").Value = 1.75
$ws.Range("W16").Value = 7.5
$ws.Range("X16").Value = 11
$ws.Range("Y16").Value = 9.5
$ws.Range("Z16").Value = 21
$ws.Range("AA16").Value = 19
$ws.Range("AH16").Value = 9.5
$ws.Range("AI16").Value = 15
$ws.Range("AK16").Value = 34
$ws.Range("AL16").Value = 26
$ws.Range("AN16").Value = 4.33
$ws.Range("AO16").Value = 13
$ws.Range("AQ16").Value = 41
$ws.Range("AS16").Value = 151
$ws.Range("AW16").Value = 5
